# Refresh the cryptocurrency market snapshot (price / 1h volume columns,
# plus two rows whose ranking happened to swap: Quant <-> Aave).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Many of the new "Price" figures are plain digits (e.g. 0.7174,
    # 243.57, 1.000). Assigning those straight to .Value lets Excel
    # auto-convert them to numbers, which is not what the sheet wants -
    # the Price column stores everything as literal text. Prefixing with
    # a single quote forces text entry (like typing it in the UI), and
    # resetting .Style afterwards drops the implicit quote-prefix style
    # so no extra formatting is left behind on the cell.
    $r = $ws.Range($cellRef)
    $r.Value = "'" + $text
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "29.386.20"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.876.62"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue "D5" "0.7174"
$ws.Range("E5").Value = "  +1.16%  "
Set-TextValue "D6" "243.57"
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("E7").Value = "  +0.06%  "
Set-TextValue "D8" "0.07953"
$ws.Range("E8").Value = "  +1.06%  "
Set-TextValue "D9" "0.3148"
$ws.Range("E9").Value = "  +1.16%  "
Set-TextValue "D10" "24.97"
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("E11").Value = "  -3.18%  "
$ws.Range("D12").Value = "1.894.16"
$ws.Range("E12").Value = "  +1.15%  "
Set-TextValue "D13" "95.59"
$ws.Range("E13").Value = "  +4.73%  "
Set-TextValue "D14" "5.235"
$ws.Range("E14").Value = "  -0.01%  "
Set-TextValue "D15" "0.7071"
$ws.Range("E15").Value = "  -1.56%  "
Set-TextValue "D16" "6.401"
$ws.Range("E16").Value = "  +4.32%  "
Set-TextValue "D17" "0.000008450"
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("D18").Value = "29.403.07"
$ws.Range("E18").Value = "  +0.31%  "
Set-TextValue "D19" "253.07"
$ws.Range("E19").Value = "  +5.07%  "
Set-TextValue "D20" "13.41"
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("D21").Value = "2.138.35"
$ws.Range("E21").Value = "  +0.61%  "
Set-TextValue "D22" "1.000"
$ws.Range("E22").Value = "  +0.06%  "
Set-TextValue "D23" "7.666"
$ws.Range("E23").Value = "  -1.20%  "
Set-TextValue "D24" "1.001"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  -0.35%  "
Set-TextValue "D26" "9.066"
$ws.Range("E26").Value = "  +0.22%  "
Set-TextValue "D27" "161.86"
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("E28").Value = "  +2.04%  "
Set-TextValue "D29" "1.509"
$ws.Range("E29").Value = "  +0.13%  "
Set-TextValue "D30" "4.419"
$ws.Range("E30").Value = "  +0.22%  "
Set-TextValue "D31" "4.296"
$ws.Range("E31").Value = "  -0.97%  "
$ws.Range("E32").Value = "  -0.07%  "
Set-TextValue "D33" "0.05326"
$ws.Range("E33").Value = "  -0.51%  "
Set-TextValue "D34" "1.945"
Set-TextValue "D35" "0.7564"
$ws.Range("E35").Value = "  +1.23%  "
$ws.Range("E36").Value = "  +0.08%  "
Set-TextValue "D38" "0.01890"
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("D39").Value = "1.267.27"
$ws.Range("E39").Value = "  +1.67%  "
Set-TextValue "D40" "2.764"
$ws.Range("E40").Value = "  +1.12%  "
Set-TextValue "D41" "6.388"
$ws.Range("E41").Value = "  -1.79%  "
Set-TextValue "D42" "0.9059"
$ws.Range("E42").Value = "  +1.62%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D43" "111.89"
$ws.Range("E43").Value = "  +1.36%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D44" "74.25"
$ws.Range("E44").Value = "  +2.48%  "
Set-TextValue "D45" "1.001"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("D47").Value = "2.034.11"
$ws.Range("E47").Value = "  +0.78%  "
Set-TextValue "D48" "1.811"
$ws.Range("E48").Value = "  +0.73%  "
Set-TextValue "D49" "0.5203"
$ws.Range("E49").Value = "  +0.32%  "
Set-TextValue "D50" "9.526"
Set-TextValue "D51" "0.4346"
$ws.Range("E51").Value = "  -0.07%  "
